$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C33").Value = "[name=`"Armed Infected`"]  You still think you’re so bright, hiding the truth to 'maintain order?'`n"
$ws.Range("C38").Value = "[name=`"Severin`"]  But, I’m sorry to say, we may be about to have an 'accident' here.`n"
$ws.Range("C40").Value = "[name=`"Severin`"]  No, just a 'slip' of the hand—`n"
$ws.Range("C51").Value = "[name=`"Folinic`"]  I’ll have words for you soon, 'Schultz.'`n"
$ws.Range("C54").Value = "[name=`"Severin`"]  I told you before, I’m not 'Sir.' I’m not your Schultz when we’re alone. I’m Thor’s father, and your uncle.`n"
